# Add a new verification test step (TS) for "lnk_Women" web element into the
# "Test Steps" sheet, inside the Verify_MyAccount_Page block (after the
# btn_Mywhishlist row), renumbering the subsequent TS_### ids, and nudge the
# two sheet-view selections, matching the author's commit.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Test Steps")
$ws2 = $wb.Worksheets.Item("Test Cases")

# --- Insert a new row 18 (pushing the old rows 18-21 down to 19-22), and
# clone the formatting of row 17 (the last "Verify_MyAccount_Page" row) onto
# it so borders/styles match the rest of the block. ---
$ws1.Rows.Item(18).Insert()
$ws1.Range("A17:G17").Copy()
$ws1.Range("A18:G18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Fill in the new row: verifying the lnk_Women element of My Account page. ---
$ws1.Range("A18").Value = "Verify_MyAccount_Page"
$ws1.Range("B18").Value = "TS_017"
$ws1.Range("C18").Value = "Verify the elements lnk_Women of My Account Page"
$ws1.Range("D18").Value = "lnk_Women"
$ws1.Range("E18").Value = "verifyElement"
$ws1.Range("F18").Value = ""
$ws1.Range("G18").Value = "PASS"

# --- Renumber the TS_### ids sequentially for every row in the
# "Verify_MyAccount_Page" block plus everything below it, since the test
# steps keep a monotonically-increasing TS id matching their physical
# position on the sheet (rows 18-21 used to be 19-22). ---
$ws1.Range("B13").Value = "TS_012"
$ws1.Range("B14").Value = "TS_013"
$ws1.Range("B15").Value = "TS_014"
$ws1.Range("B16").Value = "TS_015"
$ws1.Range("B17").Value = "TS_016"
$ws1.Range("B19").Value = "TS_018"
$ws1.Range("B20").Value = "TS_019"
$ws1.Range("B21").Value = "TS_020"
$ws1.Range("B22").Value = "TS_021"

# --- Sheet view nudges recorded in the commit. ---
$ws1.Application.ActiveWindow.ScrollRow = 1
$ws1.Range("A2").Select()
$ws1.Range("D18").Select()

$ws2.Range("A2").Select()

$wb.Save()
